$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (year 2025) metrics as per updated data
$ws.Range("C8").Value = 1033
$ws.Range("D8").Value = 170
$ws.Range("E8").Value = 863
$ws.Range("F8").Value = 6.972928630024611
$ws.Range("G8").Value = 83.5430784123911
$ws.Range("H8").Value = 16.45692158760891
